$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 7): name, gender, email, status
$ws.Range("A7").Value = "krat"
$ws.Range("B7").Value = "female"
$ws.Range("C7").Value = "krat011@gmail.com"
$ws.Range("D7").Value = "active"

# Turn the email cell into a mailto hyperlink, like the other rows
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:krat011@gmail.com")

# Hyperlinks.Add applies its own font-only style variant; put the cell back
# on the shared "Hyperlink" cell style used by the rest of column C.
$ws.Range("C7").Style = "Hyperlink"

# Update the selected/active cell to the new last row in column D
$ws.Range("D7").Select()
